# Insert a new data row at row 191 (pushing existing rows 191-280 down to
# 192-281) and populate it with the new weekly price observation for
# "Haba" at Vega Central Mapocho de Santiago, reported 2022-10-05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).Insert()

$ws.Cells.Item(191, 1).Value  = 9
$ws.Cells.Item(191, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(191, 3).Value  = "Metropolitana"
$ws.Cells.Item(191, 4).Value  = 44839
$ws.Cells.Item(191, 5).Value  = 13
$ws.Cells.Item(191, 6).Value  = 100112026
$ws.Cells.Item(191, 7).Value  = "Haba"
$ws.Cells.Item(191, 8).Value  = "Sin especificar"
$ws.Cells.Item(191, 9).Value  = "Primera"
$ws.Cells.Item(191, 10).Value = 52
$ws.Cells.Item(191, 11).Value = 11000
$ws.Cells.Item(191, 12).Value = 12000
$ws.Cells.Item(191, 13).Value = 11500
$ws.Cells.Item(191, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(191, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(191, 16).Value = 460
$ws.Cells.Item(191, 17).Value = 25
$ws.Cells.Item(191, 18).Value = "Hortaliza"
